# Remove stray "[n]" footnote markers from vaccine-name strings and
# collapse multi-line Vaccine/BrandName labels (e.g. "Hepatitis B [5]\nPediatric/Adolescent")
# into single-line text across all four sheets. Also fixes a mis-pointed
# shared string so the "Afluria Quadrivalent" brand name (instead of the stray
# "Afluria\nQuadrivalent" duplicate) is used for the Adult Influenza rows 9-10.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Pediatric VFC Vaccine ')
$ws.Range('A2').Value = 'DTaP '
$ws.Range('A3').Value = 'DTaP '
$ws.Range('A4').Value = 'DTaP-IPV '
$ws.Range('A5').Value = 'DTaP-IPV '
$ws.Range('A6').Value = 'DTaP-Hep B-IPV '
$ws.Range('A7').Value = 'DTaP-IP-HI '
$ws.Range('A8').Value = 'DTaP-IPV-HIB-HEPB '
$ws.Range('A9').Value = 'DTaP-IPV-HIB-HEPB '
$ws.Range('A10').Value = 'Dengue Tetravalent Vaccine, Live '
$ws.Range('A11').Value = 'e-IPV '
$ws.Range('A12').Value = 'Hepatitis A Pediatric '
$ws.Range('A13').Value = 'Hepatitis A Pediatric '
$ws.Range('A14').Value = 'Hepatitis A-Hepatitis B 18 only '
$ws.Range('A15').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('A16').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('B16').Value = 'Recombivax HB'
$ws.Range('A17').Value = 'Hib '
$ws.Range('A18').Value = 'Hib '
$ws.Range('A19').Value = 'Hib '
$ws.Range('A20').Value = 'HPV - Human Papillomavirus 9-valent '
$ws.Range('A21').Value = 'MENB - Meningococcal Group B '
$ws.Range('A22').Value = 'MENB - Meningococcal Group B '
$ws.Range('A23').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A24').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A25').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A26').Value = 'Measles, Mumps and Rubella (MMR) '
$ws.Range('A27').Value = 'MMR/Varicella '
$ws.Range('A28').Value = 'Pneumococcal 13-valent  (Pediatric)'
$ws.Range('A30').Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range('A31').Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range('A32').Value = 'Rotavirus, Live, Oral, Oral '
$ws.Range('A33').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A34').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A35').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A36').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A37').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A38').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A39').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A40').Value = 'Varicella '

$ws = $wb.Worksheets.Item('Adult Vaccine ')
$ws.Range('A2').Value = 'Hepatitis A Adult '
$ws.Range('A3').Value = 'Hepatitis A Adult '
$ws.Range('A4').Value = 'Hepatitis A-Hepatitis B Adult '
$ws.Range('A5').Value = 'Hepatitis B Adult '
$ws.Range('A6').Value = 'Hepatitis B Adult '
$ws.Range('A7').Value = 'Hepatitis B Adult '
$ws.Range('A8').Value = 'HPV-Human Papillomavirus 9 Valent '
$ws.Range('A9').Value = 'Measles, Mumps,  Rubella '
$ws.Range('A10').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A11').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A12').Value = 'MENB - Meningococcal Group B '
$ws.Range('A13').Value = 'MENB - Meningococcal Group B '
$ws.Range('A14').Value = 'Pneumococcal 13-valent '
$ws.Range('A16').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A17').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A18').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A19').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A20').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A21').Value = 'Varicella '

$ws = $wb.Worksheets.Item('Pediatric Influenza Vaccine ')
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B2').Value = 'Fluzone Quadrivalent'
$ws.Range('A3').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B3').Value = 'Fluzone Quadrivalent'
$ws.Range('A4').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B4').Value = 'Fluzone Quadrivalent'
$ws.Range('A5').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B5').Value = 'Fluarix Quadrivalent'
$ws.Range('A6').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B6').Value = 'FluLaval Quadrivalent'
$ws.Range('A7').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A8').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A9').Value = 'Influenza  (Age 6 -35 months)'
$ws.Range('A10').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range('B10').Value = 'FluMist Quadrivalent'

$ws = $wb.Worksheets.Item('Adult Influenza Vaccine ')
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B2').Value = 'Fluzone Quadrivalent'
$ws.Range('A3').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B3').Value = 'Fluzone Quadrivalent'
$ws.Range('A4').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B4').Value = 'Fluzone Quadrivalent'
$ws.Range('A5').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B5').Value = 'Fluarix Quadrivalent'
$ws.Range('A6').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B6').Value = 'FluLaval Quadrivalent'
$ws.Range('A7').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A8').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A9').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('B9').Value = 'Afluria Quadrivalent'
$ws.Range('A10').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B10').Value = 'Afluria Quadrivalent'
$ws.Range('A11').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range('B11').Value = 'FluMist Quadrivalent'
